$wb = $excel.ActiveWorkbook

# Fold_1
$ws = $wb.Worksheets.Item("Fold_1")
$ws.Cells.Item(2, 2).Value = 6.87565
$ws.Cells.Item(2, 3).Value = 1.2232
$ws.Cells.Item(2, 4).Value = 35.04514999999999
$ws.Cells.Item(2, 5).Value = 0.0359
$ws.Cells.Item(2, 6).Value = 26.92725
$ws.Cells.Item(2, 7).Value = 26.92725
$ws.Cells.Item(2, 11).Value = 67.2974
$ws.Cells.Item(2, 12).Value = 26.927
$ws.Cells.Item(2, 13).Value = 40.3704
$ws.Cells.Item(2, 14).Value = 7.785999999999999
$ws.Cells.Item(2, 15).Value = 32.5844
$ws.Cells.Item(3, 2).Value = 12.049
$ws.Cells.Item(3, 3).Value = 3.54
$ws.Cells.Item(3, 4).Value = 53.524
$ws.Cells.Item(3, 6).Value = 30.926
$ws.Cells.Item(3, 7).Value = 29.193
$ws.Cells.Item(3, 8).Value = 1.733
$ws.Cells.Item(3, 9).Value = 1.733
$ws.Cells.Item(3, 11).Value = 32.437
$ws.Cells.Item(3, 12).Value = 29.193
$ws.Cells.Item(3, 13).Value = 3.244
$ws.Cells.Item(3, 14).Value = 3.244
$ws.Cells.Item(4, 2).Value = 33.529
$ws.Cells.Item(4, 3).Value = 2.449
$ws.Cells.Item(4, 4).Value = 53.46
$ws.Cells.Item(4, 6).Value = 35.893
$ws.Cells.Item(4, 7).Value = 35.8925808219178
$ws.Cells.Item(4, 11).Value = 35.893
$ws.Cells.Item(4, 12).Value = 35.893

# Fold_2
$ws = $wb.Worksheets.Item("Fold_2")
$ws.Cells.Item(2, 2).Value = 6.19875
$ws.Cells.Item(2, 3).Value = 1.20455
$ws.Cells.Item(2, 4).Value = 34.63824999999999
$ws.Cells.Item(2, 5).Value = 0.0359
$ws.Cells.Item(2, 6).Value = 26.70205
$ws.Cells.Item(2, 7).Value = 26.70205
$ws.Cells.Item(2, 11).Value = 229.8882
$ws.Cells.Item(2, 12).Value = 26.70200000000001
$ws.Cells.Item(2, 13).Value = 203.1862
$ws.Cells.Item(2, 14).Value = 8.9472
$ws.Cells.Item(2, 15).Value = 194.2392
$ws.Cells.Item(3, 2).Value = 9.163
$ws.Cells.Item(3, 3).Value = 3.818
$ws.Cells.Item(3, 4).Value = 53.394
$ws.Cells.Item(3, 6).Value = 29.79
$ws.Cells.Item(3, 7).Value = 28.304
$ws.Cells.Item(3, 8).Value = 1.486
$ws.Cells.Item(3, 9).Value = 1.486
$ws.Cells.Item(3, 11).Value = 109.947
$ws.Cells.Item(3, 12).Value = 28.30399999999999
$ws.Cells.Item(3, 13).Value = 81.643
$ws.Cells.Item(3, 14).Value = 7.8102
$ws.Cells.Item(3, 15).Value = 73.83279999999999
$ws.Cells.Item(4, 2).Value = 19.899
$ws.Cells.Item(4, 3).Value = 2.786
$ws.Cells.Item(4, 4).Value = 53.876
$ws.Cells.Item(4, 6).Value = 31.612
$ws.Cells.Item(4, 7).Value = 31.61164383561644
$ws.Cells.Item(4, 11).Value = 35.071
$ws.Cells.Item(4, 12).Value = 31.612
$ws.Cells.Item(4, 13).Value = 3.459000000000001
$ws.Cells.Item(4, 14).Value = 3.459000000000001

# Fold_3
$ws = $wb.Worksheets.Item("Fold_3")
$ws.Cells.Item(2, 2).Value = 7.723350000000001
$ws.Cells.Item(2, 3).Value = 1.3063
$ws.Cells.Item(2, 4).Value = 32.7205
$ws.Cells.Item(2, 6).Value = 27.15645
$ws.Cells.Item(2, 7).Value = 27.15645
$ws.Cells.Item(2, 11).Value = 77.2728
$ws.Cells.Item(2, 12).Value = 27.156
$ws.Cells.Item(2, 13).Value = 50.1168
$ws.Cells.Item(2, 14).Value = 5.8988
$ws.Cells.Item(2, 15).Value = 44.218
$ws.Cells.Item(3, 2).Value = 12.049
$ws.Cells.Item(3, 3).Value = 3.54
$ws.Cells.Item(3, 4).Value = 53.524
$ws.Cells.Item(3, 6).Value = 31.738
$ws.Cells.Item(3, 7).Value = 29.193
$ws.Cells.Item(3, 8).Value = 2.544
$ws.Cells.Item(3, 9).Value = 2.544
$ws.Cells.Item(3, 11).Value = 29.193
$ws.Cells.Item(3, 12).Value = 29.193
$ws.Cells.Item(4, 2).Value = 33.529
$ws.Cells.Item(4, 3).Value = 2.449
$ws.Cells.Item(4, 4).Value = 53.46
$ws.Cells.Item(4, 6).Value = 35.893
$ws.Cells.Item(4, 7).Value = 35.8925808219178
$ws.Cells.Item(4, 11).Value = 35.893
$ws.Cells.Item(4, 12).Value = 35.893

# Fold_4
$ws = $wb.Worksheets.Item("Fold_4")
$ws.Cells.Item(2, 2).Value = 8.298950000000001
$ws.Cells.Item(2, 3).Value = 1.2959
$ws.Cells.Item(2, 4).Value = 32.83895
$ws.Cells.Item(2, 5).Value = 0.0359
$ws.Cells.Item(2, 6).Value = 27.3412
$ws.Cells.Item(2, 7).Value = 27.3412
$ws.Cells.Item(2, 11).Value = 47.029
$ws.Cells.Item(2, 12).Value = 27.341
$ws.Cells.Item(2, 13).Value = 19.688
$ws.Cells.Item(2, 14).Value = 2.7378
$ws.Cells.Item(2, 15).Value = 16.9502
$ws.Cells.Item(3, 2).Value = 12.049
$ws.Cells.Item(3, 3).Value = 3.54
$ws.Cells.Item(3, 4).Value = 53.524
$ws.Cells.Item(3, 6).Value = 31.738
$ws.Cells.Item(3, 7).Value = 29.193
$ws.Cells.Item(3, 8).Value = 2.544
$ws.Cells.Item(3, 9).Value = 2.544
$ws.Cells.Item(3, 11).Value = 29.193
$ws.Cells.Item(3, 12).Value = 29.193
$ws.Cells.Item(4, 2).Value = 33.529
$ws.Cells.Item(4, 3).Value = 2.449
$ws.Cells.Item(4, 4).Value = 53.46
$ws.Cells.Item(4, 6).Value = 35.893
$ws.Cells.Item(4, 7).Value = 35.8925808219178
$ws.Cells.Item(4, 11).Value = 35.893
$ws.Cells.Item(4, 12).Value = 35.893

# Fold_5
$ws = $wb.Worksheets.Item("Fold_5")
$ws.Cells.Item(2, 2).Value = 8.0367
$ws.Cells.Item(2, 3).Value = 1.05765
$ws.Cells.Item(2, 4).Value = 33.68215
$ws.Cells.Item(2, 5).Value = 0.0359
$ws.Cells.Item(2, 6).Value = 27.25005000000001
$ws.Cells.Item(2, 7).Value = 27.25005000000001
$ws.Cells.Item(2, 11).Value = 145.9106
$ws.Cells.Item(2, 12).Value = 27.25
$ws.Cells.Item(2, 13).Value = 118.6606
$ws.Cells.Item(2, 14).Value = 3.2118
$ws.Cells.Item(2, 15).Value = 115.4488
$ws.Cells.Item(3, 2).Value = 12.081
$ws.Cells.Item(3, 3).Value = 2.901
$ws.Cells.Item(3, 4).Value = 52.019
$ws.Cells.Item(3, 6).Value = 31.645
$ws.Cells.Item(3, 7).Value = 29.105
$ws.Cells.Item(3, 8).Value = 2.54
$ws.Cells.Item(3, 9).Value = 2.54
$ws.Cells.Item(3, 11).Value = 42.3792
$ws.Cells.Item(3, 12).Value = 29.105
$ws.Cells.Item(3, 13).Value = 13.2742
$ws.Cells.Item(3, 14).Value = 0.2958
$ws.Cells.Item(3, 15).Value = 12.9786
$ws.Cells.Item(4, 2).Value = 33.529
$ws.Cells.Item(4, 3).Value = 2.449
$ws.Cells.Item(4, 4).Value = 52.947
$ws.Cells.Item(4, 6).Value = 35.882
$ws.Cells.Item(4, 7).Value = 35.88203287671233
$ws.Cells.Item(4, 11).Value = 35.926
$ws.Cells.Item(4, 12).Value = 35.882
$ws.Cells.Item(4, 13).Value = 0.044
$ws.Cells.Item(4, 14).Value = 0.044
$ws.Cells.Item(4, 15).Value = 0
